$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set L3 to 2021 (was 2020)
$ws.Range("L3").Value = 2021

# Add new column M, row 3 (2022) - copy L3's format first, then set the value
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("M3").Value = 2022

# Add new column M, row 4 (6.18) - copy L4's format first, then set the value
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 6.18

$excel.CutCopyMode = $false

# Update selection to M9 (was M12)
$ws.Range("M9").Select()
